# Applies the updated pricing/profit figures captured in the diff.
# Values originate from a scheduled pricing-data refresh (Universalis API snapshot).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 2 (Leve Item ID 5489)
$ws.Range("H2").Value = 1082.8
$ws.Range("I2").Value = 92.111115
$ws.Range("J2").Value = 9999
$ws.Range("K2").Value = 92.111115
$ws.Range("L2").Value = 9999
$ws.Range("M2").Value = 20.888885
$ws.Range("N2").Value = -10225

# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 95.30768999999999
$ws.Range("I5").Value = 95.30768999999999
$ws.Range("K5").Value = 95.30768999999999
$ws.Range("M5").Value = 19.69231000000001

# Row 98 (Leve Item ID 36237)
$ws.Range("H98").Value = 1332
$ws.Range("I98").Value = 1332
$ws.Range("K98").Value = 1332
$ws.Range("M98").Value = 166

# Row 103 (Leve Item ID 19909)
$ws.Range("H103").Value = 5226.5
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 116 (Leve Item ID 27778)
$ws.Range("H116").Value = 3465.889
$ws.Range("I116").Value = 3198.8
$ws.Range("J116").Value = 3799.75
$ws.Range("K116").Value = 3198.8
$ws.Range("L116").Value = 3799.75
$ws.Range("M116").Value = 243.1999999999998
$ws.Range("N116").Value = -10683.75

# Row 122 (Leve Item ID 36237)
$ws.Range("H122").Value = 1332
$ws.Range("I122").Value = 1332
$ws.Range("K122").Value = 3996
$ws.Range("M122").Value = -1546

# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2711.5122
$ws.Range("I138").Value = 1696.6364
$ws.Range("J138").Value = 3083.6333
$ws.Range("K138").Value = 5089.9092
$ws.Range("L138").Value = 9250.8999
$ws.Range("M138").Value = 50.09079999999994
$ws.Range("N138").Value = -19530.8999

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 4434.8936
$ws.Range("I32").Value = 4434.8936
$ws.Range("K32").Value = 4434.8936
$ws.Range("M32").Value = -4147.8936

# Row 44 (Leve Item ID 3861)
$ws.Range("H44").Value = 25500
$ws.Range("J44").Value = 31000
$ws.Range("L44").Value = 31000
$ws.Range("N44").Value = -31976

# Row 55 (Leve Item ID 2830)
$ws.Range("H55").Value = 19950
$ws.Range("J55").Value = 19900
$ws.Range("L55").Value = 19900
$ws.Range("N55").Value = -20530

# Row 110 (Leve Item ID 27708)
$ws.Range("H110").Value = 6304.9414
$ws.Range("I110").Value = 4461.875
$ws.Range("J110").Value = 7943.222
$ws.Range("K110").Value = 4461.875
$ws.Range("L110").Value = 7943.222
$ws.Range("M110").Value = -2416.875
$ws.Range("N110").Value = -12033.222

# Row 122 (Leve Item ID 36168)
$ws.Range("H122").Value = 1449.9445
$ws.Range("I122").Value = 1447
$ws.Range("K122").Value = 4341
$ws.Range("M122").Value = -1891

$ws = $wb.Worksheets.Item("BSM")
# Row 55 (Leve Item ID 27151)
$ws.Range("H55").Value = 73293
$ws.Range("J55").Value = 73293
$ws.Range("L55").Value = 73293
$ws.Range("N55").Value = -73839

# Row 94 (Leve Item ID 19939)
$ws.Range("H94").Value = 2436
$ws.Range("I94").Value = 1973.4286
$ws.Range("K94").Value = 1973.4286
$ws.Range("M94").Value = -1522.4286

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 29416774
$ws.Range("I31").Value = 52634396
$ws.Range("K31").Value = 52634396
$ws.Range("M31").Value = -52634101

# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 29416774
$ws.Range("I34").Value = 52634396
$ws.Range("K34").Value = 52634396
$ws.Range("M34").Value = -52634194

# Row 41 (Leve Item ID 1917)
$ws.Range("H41").Value = 21244.5
$ws.Range("J41").Value = 21244.5
$ws.Range("L41").Value = 21244.5
$ws.Range("N41").Value = -22100.5

# Row 59 (Leve Item ID 1942)
$ws.Range("H59").Value = 22612.666
$ws.Range("J59").Value = 22612.666
$ws.Range("L59").Value = 22612.666
$ws.Range("N59").Value = -24902.666

# Row 60 (Leve Item ID 1937)
$ws.Range("H60").Value = 39200
$ws.Range("I60").Value = 70000
$ws.Range("J60").Value = 23800
$ws.Range("K60").Value = 70000
$ws.Range("L60").Value = 23800
$ws.Range("M60").Value = -69489
$ws.Range("N60").Value = -24822

# Row 64 (Leve Item ID 10610)
$ws.Range("H64").Value = 49970
$ws.Range("J64").Value = 49970
$ws.Range("L64").Value = 49970
$ws.Range("N64").Value = -50466

# Row 67 (Leve Item ID 10610)
$ws.Range("H67").Value = 49970
$ws.Range("J67").Value = 49970
$ws.Range("L67").Value = 49970
$ws.Range("N67").Value = -51686

# Row 68 (Leve Item ID 10611)
$ws.Range("H68").Value = 40559.4
$ws.Range("J68").Value = 39999.25
$ws.Range("L68").Value = 39999.25
$ws.Range("N68").Value = -41497.25

# Row 71 (Leve Item ID 10611)
$ws.Range("H71").Value = 40559.4
$ws.Range("J71").Value = 39999.25
$ws.Range("L71").Value = 119997.75
$ws.Range("N71").Value = -127485.75

# Row 122 (Leve Item ID 36196)
$ws.Range("H122").Value = 3252.111
$ws.Range("I122").Value = 3475.8667
$ws.Range("J122").Value = 2133.3333
$ws.Range("K122").Value = 10427.6001
$ws.Range("L122").Value = 6399.999899999999
$ws.Range("M122").Value = -7977.6001
$ws.Range("N122").Value = -11299.9999

# Row 132 (Leve Item ID 44019)
$ws.Range("H132").Value = 2909.2
$ws.Range("I132").Value = 2212.7058
$ws.Range("K132").Value = 6638.117400000001
$ws.Range("M132").Value = -4108.117400000001

$ws = $wb.Worksheets.Item("CUL")
# Row 68 (Leve Item ID 12895)
$ws.Range("H68").Value = 50003030
$ws.Range("I68").Value = 100002424
$ws.Range("J68").Value = 3639.4
$ws.Range("K68").Value = 300007272
$ws.Range("L68").Value = 10918.2
$ws.Range("M68").Value = -300006461
$ws.Range("N68").Value = -12540.2

# Row 71 (Leve Item ID 12895)
$ws.Range("H71").Value = 50003030
$ws.Range("I71").Value = 100002424
$ws.Range("J71").Value = 3639.4
$ws.Range("K71").Value = 900021816
$ws.Range("L71").Value = 32754.6
$ws.Range("M71").Value = -900017760
$ws.Range("N71").Value = -40866.60000000001

# Row 92 (Leve Item ID 19841)
$ws.Range("H92").Value = 448.25
$ws.Range("J92").Value = 448.25
$ws.Range("L92").Value = 1344.75
$ws.Range("N92").Value = -3840.75

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (Leve Item ID 5062)
$ws.Range("H2").Value = 69.69231000000001
$ws.Range("I2").Value = 63.5
$ws.Range("K2").Value = 63.5
$ws.Range("M2").Value = 49.5

# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 22850
$ws.Range("J46").Value = 25239
$ws.Range("L46").Value = 25239
$ws.Range("N46").Value = -25551

# Row 80 (Leve Item ID 12521)
$ws.Range("H80").Value = 4494.2
$ws.Range("I80").Value = 4390.5835
$ws.Range("J80").Value = 4649.625
$ws.Range("K80").Value = 4390.5835
$ws.Range("L80").Value = 4649.625
$ws.Range("M80").Value = -3392.5835
$ws.Range("N80").Value = -6645.625

# Row 83 (Leve Item ID 12521)
$ws.Range("H83").Value = 4494.2
$ws.Range("I83").Value = 4390.5835
$ws.Range("J83").Value = 4649.625
$ws.Range("K83").Value = 21952.9175
$ws.Range("L83").Value = 23248.125
$ws.Range("M83").Value = -16960.9175
$ws.Range("N83").Value = -33232.125

# Row 113 (Leve Item ID 27710)
$ws.Range("H113").Value = 1592.125
$ws.Range("I113").Value = 1185.25
$ws.Range("J113").Value = 1999
$ws.Range("K113").Value = 1185.25
$ws.Range("L113").Value = 1999
$ws.Range("M113").Value = 984.75
$ws.Range("N113").Value = -6339

# Row 116 (Leve Item ID 26120)
$ws.Range("H116").Value = 0
$ws.Range("J116").Value = 0
$ws.Range("L116").Value = 0
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
# Row 82 (Leve Item ID 12565)
$ws.Range("H82").Value = 1913.2858
$ws.Range("I82").Value = 1932.6666
$ws.Range("J82").Value = 1797
$ws.Range("K82").Value = 1932.6666
$ws.Range("L82").Value = 1797
$ws.Range("M82").Value = -1571.6666
$ws.Range("N82").Value = -2519

# Row 85 (Leve Item ID 12565)
$ws.Range("H85").Value = 1913.2858
$ws.Range("I85").Value = 1932.6666
$ws.Range("J85").Value = 1797
$ws.Range("K85").Value = 1932.6666
$ws.Range("L85").Value = 1797
$ws.Range("M85").Value = -684.6666
$ws.Range("N85").Value = -4293

# Row 122 (Leve Item ID 36247)
$ws.Range("H122").Value = 5618.923
$ws.Range("I122").Value = 3607.9
$ws.Range("J122").Value = 12322.333
$ws.Range("K122").Value = 10823.7
$ws.Range("L122").Value = 36966.999
$ws.Range("M122").Value = -8373.700000000001
$ws.Range("N122").Value = -41866.999

# Row 132 (Leve Item ID 44058)
$ws.Range("H132").Value = 11422.037
$ws.Range("I132").Value = 17519.666
$ws.Range("J132").Value = 3800
$ws.Range("K132").Value = 52558.99800000001
$ws.Range("L132").Value = 11400
$ws.Range("M132").Value = -50028.99800000001
$ws.Range("N132").Value = -16460

$ws = $wb.Worksheets.Item("WVR")
# Row 51 (Leve Item ID 3162)
$ws.Range("H51").Value = 21035
$ws.Range("J51").Value = 32000
$ws.Range("L51").Value = 32000
$ws.Range("N51").Value = -33020

# Row 111 (Leve Item ID 25833)
$ws.Range("H111").Value = 88777
$ws.Range("J111").Value = 88777
$ws.Range("L111").Value = 88777
$ws.Range("N111").Value = -96957

# Row 122 (Leve Item ID 36208)
$ws.Range("H122").Value = 3169.8809
$ws.Range("I122").Value = 2336.742
$ws.Range("K122").Value = 7010.226000000001
$ws.Range("M122").Value = -4560.226000000001

# Row 124 (Leve Item ID 34280)
$ws.Range("H124").Value = 100000
$ws.Range("J124").Value = 100000
$ws.Range("L124").Value = 100000
$ws.Range("N124").Value = -109820

# Row 132 (Leve Item ID 44029)
$ws.Range("H132").Value = 4337.9756
$ws.Range("I132").Value = 3633.625
$ws.Range("J132").Value = 6842.3335
$ws.Range("K132").Value = 10900.875
$ws.Range("L132").Value = 20527.0005
$ws.Range("M132").Value = -8370.875
$ws.Range("N132").Value = -25587.0005

# Row 136 (Leve Item ID 44031)
$ws.Range("H136").Value = 1655.3334
$ws.Range("I136").Value = 1190
$ws.Range("J136").Value = 5098.8
$ws.Range("K136").Value = 3570
$ws.Range("L136").Value = 15296.4
$ws.Range("M136").Value = -1020
$ws.Range("N136").Value = -20396.4
